$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Motor pulley row (row 6): the link text now points to the
# alternate AliExpress listing that was previously only mentioned in the
# note, and that note is no longer needed now that it's the chosen part.
$ws.Range("B6").Value = "https://www.aliexpress.us/item/3256805934638435.html"
$ws.Range("F6").Clear()

# The note that made this row tall is gone, so let the row shrink back down.
$ws.Rows("6:6").AutoFit()

# Update the currently selected cell to match the saved view state.
$ws.Range("G7").Select()
